# Exchange rate.xlsx - regenerate with updated NBP data:
#  - two extra trading days (2021-03-18, 2021-03-19) appended to every
#    currency sheet
#  - the four sheets were rebuilt/reshuffled by the source tool: sheet1
#    now carries the EUR table, sheet2 the USD table (they swapped),
#    sheet3 (GBP) stays put, and sheet4 switched from CNY to CAD
#  - every sheet tab name got a new trailing "id <random>" suffix

function Set-TextValue($ws, $addr, $text) {
    # Assign $text as a literal (shared-string) value, never letting Excel's
    # type-sniffer reinterpret it as a number/date. We build it via a
    # formula that evaluates to the literal text, then paste-special just
    # the computed value back over itself (drops the formula, keeps the
    # cell's existing style untouched).
    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

function Add-RateRow($ws, $rowNum, $templateRowNum, $no, $date, $mid) {
    # Write a brand-new data row (no / effectiveDate / mid in D:F), then
    # copy the formatting from an existing data row so the new row matches
    # (same style index, no stray styles created).
    Set-TextValue $ws ("D" + $rowNum) $no
    Set-TextValue $ws ("E" + $rowNum) $date
    Set-TextValue $ws ("F" + $rowNum) $mid

    $ws.Range("D" + $templateRowNum + ":F" + $templateRowNum).Copy()
    $ws.Range("D" + $rowNum + ":F" + $rowNum).PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- sheet1: was "dolar amerykanski", becomes "euro" -----------------
Set-TextValue $ws1 "B2" "euro"
Set-TextValue $ws1 "C2" "EUR"
$eurRates = @("4.5231","4.5345","4.5393","4.554","4.5793","4.5903","4.5844","4.5718","4.5805","4.5909","4.5836","4.5949","4.6065")
for ($i = 0; $i -lt $eurRates.Length; $i++) {
    Set-TextValue $ws1 ("F" + (3 + $i)) $eurRates[$i]
}
Add-RateRow $ws1 16 15 "053/A/NBP/2021" "2021-03-18" "4.6224"
Add-RateRow $ws1 17 16 "054/A/NBP/2021" "2021-03-19" "4.6226"
$ws1.Name = "euro id 0.38009425748389225"

# --- sheet2: was "euro", becomes "dolar amerykanski" ------------------
Set-TextValue $ws2 "B2" "dolar amerykański"
Set-TextValue $ws2 "C2" "USD"
$usdRates = @("3.7572","3.7765","3.7509","3.7851","3.8393","3.8665","3.8507","3.842","3.8287","3.8521","3.8429","3.8519","3.8676")
for ($i = 0; $i -lt $usdRates.Length; $i++) {
    Set-TextValue $ws2 ("F" + (3 + $i)) $usdRates[$i]
}
Add-RateRow $ws2 16 15 "053/A/NBP/2021" "2021-03-18" "3.8705"
Add-RateRow $ws2 17 16 "054/A/NBP/2021" "2021-03-19" "3.8865"
$ws2.Name = "dolar amerykański id 0.92024677"

# --- sheet3: stays "funt szterling" (GBP), values unchanged -----------
Add-RateRow $ws3 16 15 "053/A/NBP/2021" "2021-03-18" "5.4038"
Add-RateRow $ws3 17 16 "054/A/NBP/2021" "2021-03-19" "5.4116"
$ws3.Name = "funt szterling id 0.40021683445"

# --- sheet4: was "yuan renminbi (Chiny)" (CNY), becomes "dolar kanadyjski" (CAD)
Set-TextValue $ws4 "B2" "dolar kanadyjski"
Set-TextValue $ws4 "C2" "CAD"
$cadRates = @("2.9589","2.9803","2.9771","2.9906","3.0231","3.0477","3.0533","3.0374","3.0434","3.0669","3.0861","3.0848","3.1017")
for ($i = 0; $i -lt $cadRates.Length; $i++) {
    Set-TextValue $ws4 ("F" + (3 + $i)) $cadRates[$i]
}
Add-RateRow $ws4 16 15 "053/A/NBP/2021" "2021-03-18" "3.1206"
Add-RateRow $ws4 17 16 "054/A/NBP/2021" "2021-03-19" "3.1125"
$ws4.Name = "dolar kanadyjski id 0.140940833"

Write-Output "edit applied"
